$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(145, 3).Value = 1.168885292845949
$ws.Cells.Item(147, 3).Value = 1.21486052238984
$ws.Cells.Item(148, 3).Value = 1.279036382914401
$ws.Cells.Item(152, 3).Value = 1.236072833604953
$ws.Cells.Item(154, 3).Value = 1.22415486547171
$ws.Cells.Item(155, 3).Value = 1.249077098369515
$ws.Cells.Item(159, 3).Value = 1.224516623136881
$ws.Cells.Item(161, 3).Value = 1.235914030567496
$ws.Cells.Item(162, 3).Value = 1.208582765258972
$ws.Cells.Item(164, 3).Value = 1.223339205264765
$ws.Cells.Item(165, 3).Value = 1.224129108051137
$ws.Cells.Item(166, 3).Value = 1.235573301768256
$ws.Cells.Item(175, 3).Value = 1.261125837412918
$ws.Cells.Item(176, 3).Value = 1.249623747082607
$ws.Cells.Item(178, 3).Value = 1.24948293342525
$ws.Cells.Item(179, 3).Value = 1.24989081645374
$ws.Cells.Item(184, 3).Value = 1.241279612420844
$ws.Cells.Item(185, 3).Value = 1.240343896002161
$ws.Cells.Item(186, 3).Value = 1.239565513799398
$ws.Cells.Item(187, 3).Value = 1.226372595977206
$ws.Cells.Item(188, 3).Value = 1.22690209228363
$ws.Cells.Item(189, 3).Value = 1.22779570041693
$ws.Cells.Item(190, 3).Value = 1.205641576628107
$ws.Cells.Item(191, 3).Value = 1.216840934226896
$ws.Cells.Item(192, 3).Value = 1.227545350993828
$ws.Cells.Item(193, 3).Value = 1.226995105162298
$ws.Cells.Item(194, 3).Value = 1.228583304476012
$ws.Cells.Item(195, 3).Value = 1.202035027553715
$ws.Cells.Item(196, 3).Value = 1.201453220846558
$ws.Cells.Item(197, 3).Value = 1.200773896767373
$ws.Cells.Item(201, 3).Value = 1.149225999358948
$ws.Cells.Item(202, 3).Value = 1.173906296295284
$ws.Cells.Item(207, 3).Value = 1.170979371306329
$ws.Cells.Item(208, 3).Value = 1.170410974266205
$ws.Cells.Item(209, 3).Value = 1.211399886366022
$ws.Cells.Item(210, 3).Value = 1.209062815901802
$ws.Cells.Item(211, 3).Value = 1.209767756797069
$ws.Cells.Item(212, 3).Value = 1.208126048465026
$ws.Cells.Item(213, 3).Value = 1.206143856762223
$ws.Cells.Item(214, 3).Value = 1.219685195887786
$ws.Cells.Item(215, 3).Value = 1.2294443467764
$ws.Cells.Item(216, 3).Value = 1.252166326270286
$ws.Cells.Item(217, 3).Value = 1.263445535777546
$ws.Cells.Item(218, 3).Value = 1.264352006479259
$ws.Cells.Item(219, 3).Value = 1.264027669521014
$ws.Cells.Item(220, 3).Value = 1.263709736548671
$ws.Cells.Item(221, 3).Value = 1.263976970690259
$ws.Cells.Item(222, 3).Value = 1.26387563390641
$ws.Cells.Item(223, 3).Value = 1.278942049283157
$ws.Cells.Item(224, 3).Value = 1.275572667932509
$ws.Cells.Item(225, 3).Value = 1.225329112322245
$ws.Cells.Item(226, 3).Value = 1.22508264222341
$ws.Cells.Item(227, 3).Value = 1.224625767606213
$ws.Cells.Item(228, 3).Value = 1.235961747676092
$ws.Cells.Item(252, 3).Value = 1.505024410882924
$ws.Cells.Item(253, 3).Value = 1.504056378694462
$ws.Cells.Item(256, 3).Value = 1.521293420386913
$ws.Cells.Item(258, 3).Value = 1.564888318697708
$ws.Cells.Item(259, 3).Value = 1.565272833530306
$ws.Cells.Item(260, 3).Value = 1.511848028050094
$ws.Cells.Item(261, 3).Value = 1.518667350800931
$ws.Cells.Item(262, 3).Value = 1.517943032579118
$ws.Cells.Item(263, 3).Value = 1.517641299295436
